$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.160040855407715
$ws.Range("B1").Value = 2.408264875411987
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.382341384887695
$ws.Range("E1").Value = 1.227702856063843
